# lineEstimateTestData.xlsx - "Added some files to line estimate feature"
#
# Update the sample identifier values used by the line-estimate functional
# test fixtures (work / admin-sanction / technical-sanction reference
# numbers), and move the active sheet/selection from adminSanctionDetails
# to technicalSanctionDetails.

$wb = $excel.ActiveWorkbook

# workDetails sheet: abstractEstimateNumber (E2) and workIdentificationNumber (G2)
$wsWork = $wb.Worksheets.Item("workDetails")
$wsWork.Range("E2").Value = "AEN765"
$wsWork.Range("G2").Value = "WIN765"
[void]$wsWork.Range("G2").Select()

# adminSanctionDetails sheet: administrativeSanctionNumber (B2)
$wsAdmin = $wb.Worksheets.Item("adminSanctionDetails")
$wsAdmin.Range("B2").Value = "ADN765"
[void]$wsAdmin.Range("B2").Select()

# technicalSanctionDetails sheet: technicalSanctionNumber (B2).
# Activate this sheet last so it ends up the workbook's active/selected
# tab (matches activeTab moving from 4 to 5, and the tabSelected flag
# flipping from adminSanctionDetails to technicalSanctionDetails).
$wsTech = $wb.Worksheets.Item("technicalSanctionDetails")
$wsTech.Range("B2").Value = "TSN765"
$wsTech.Activate()
[void]$wsTech.Range("B2").Select()
